# DPLKINV001-019 - update Bank IDs (PAU0269 -> PAU0271 / PAU0268), trim trailing
# semicolons on long multi-line notes, and refresh the sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Tambah Setup Profile Bank" ---------------------------------
$ws.Range("F2").Value = "Username : Putri;`nPassword : bni1234/;`nKode Buku : BPU0223;`nValuta : USD - United States Dollar;`nNo Rekening : 12345678;`nPemilik Rekening : Tester;`nBank ID : -;`nEmiten : E0001 - Pemerintah;`nKode BI : 002 - Bank Rakyat Indonesia;`nCabang Bank : Pejompongan;`nAlamat : Permata Hijau II;`nNo. Telp : 999999999;`nNama PIC : Tester;`nDays Basis : 365;`nJenis Perhitungan Bunga : Memperhitungkan Hari Libur;`nSyariah : Checked;`nTgl. Non Aktif : 18/11/2022;`nNon. Aktif : Checked"
$ws.Range("M2").Value = "PAU0271"

# --- Row 3: "View Setup Profile Bank" -----------------------------------
$ws.Range("F3").Value = "Username : Putri;`nPassword : bni1234/;`nBank ID : PAU0271"
$ws.Range("M3").Value = "PAU0271"

# --- Row 4: "Verifikasi Data Setup Profile Bank" (data dikembalikan) ---
$ws.Range("F4").Value = "Username : Putri;`nPassword : bni1234/;`nBank ID : PAU0271 (sesuaikan dengan hasil generate);`nPetugas Submit : Putri;`nTanggal Verifikasi : Tanggal hari ini;`nStatus Verifikasi : 0 - Dikembalikan ke Data Entry;`nKeterangan Verifikasi : Tolong Diperbaiki"
$ws.Range("M4").Value = "PAU0271"

# --- Row 5: "Ubah Setup Profile Bank" -----------------------------------
$ws.Range("M5").Value = "PAU0271"

# --- Row 6: "Verifikasi Data Setup Profile Bank" (data disetujui) ------
$ws.Range("F6").Value = "Username : Putri;`nPassword : bni1234/;`nBank ID : PAU0271 (sesuaikan dengan hasil generate);`nPetugas Submit : Putri;`nTanggal Verifikasi : Tanggal hari ini;`nStatus Verifikasi : 1 - Setuju;`nKeterangan Verifikasi w: Disetujui"
$ws.Range("M6").Value = "PAU0271"

# --- Row 7: "Hapus Setup Profile Bank" ----------------------------------
$ws.Range("F7").Value = "Username : Putri;`nPassword : bni1234/;`nBank ID : PAU0268"

# --- Refresh sheet view (scroll position, zoom, selection) -------------
$win = $excel.ActiveWindow
$win.Zoom = 70
$win.ScrollRow = 3
$win.ScrollColumn = 2
$win.TopLeftCell = $ws.Range("B3")
$ws.Range("F7").Select()
